$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency list refreshed with new Price (D) / Volume(1h) (E) figures.
# Coin name (B) and Link (C) stay the same for these rows; only D and/or E change.
$priceVolumeUpdates = @(
    @{Row=2; D="88.738.96"; E="  +9.08%  "},
    @{Row=3; D="3.325.82"; E="  +3.75%  "},
    @{Row=4; D="1.00"; E="  +0.22%  "},
    @{Row=5; D="219.63"; E="  +4.03%  "},
    @{Row=6; D="657.04"; E="  +2.52%  "},
    @{Row=7; D="0.358"; E="  +21.28%  "},
    @{Row=8; D="0.999"; E="  +0.06%  "},
    @{Row=9; D="0.606"; E="  +2.10%  "},
    @{Row=10; D="3.325.35"; E="  +3.93%  "},
    @{Row=11; D="0.586"; E="  -1.86%  "},
    @{Row=12; D="0.0000271"; E="  +0.32%  "},
    @{Row=13; D="35.84"; E="  +11.00%  "},
    @{Row=14; D=$null; E="  +1.56%  "},
    @{Row=15; D="3.940.22"; E="  +3.93%  "},
    @{Row=16; D="5.49"; E="  +2.29%  "},
    @{Row=17; D="88.568.64"; E="  +9.17%  "},
    @{Row=18; D="3.327.58"; E="  +3.84%  "},
    @{Row=19; D="14.75"; E="  +1.87%  "},
    @{Row=20; D=$null; E="  +0.10%  "},
    @{Row=21; D="460.97"; E="  +2.49%  "},
    @{Row=22; D="9.25"; E="  -0.81%  "},
    @{Row=23; D="5.54"; E="  +5.05%  "},
    @{Row=24; D="5.59"; E="  +11.18%  "},
    @{Row=25; D="12.67"; E="  +12.82%  "},
    @{Row=26; D="3.506.77"; E="  +4.09%  "},
    @{Row=27; D="78.88"; E="  +1.57%  "},
    @{Row=28; D="0.212"; E="  +68.47%  "},
    @{Row=29; D=$null; E="  +0.65%  "},
    @{Row=30; D=$null; E="  -0.17%  "},
    @{Row=31; D="613.79"; E="  +8.54%  "},
    @{Row=32; D="9.43"; E="  +1.98%  "},
    @{Row=33; D="1.62"; E="  +8.48%  "},
    @{Row=34; D=$null; E="  +0.37%  "},
    @{Row=35; D="2.09"; E="  +2.12%  "},
    @{Row=38; D="23.69"; E="  +2.15%  "},
    @{Row=39; D="2.20"; E="  +4.53%  "},
    @{Row=40; D="0.421"; E="  +1.83%  "},
    @{Row=41; D="21.85"; E="  +5.00%  "},
    @{Row=42; D="0.999"; E="  +0.18%  "},
    @{Row=43; D="3.02"; E="  +5.18%  "},
    @{Row=44; D=$null; E="  +0.00%  "},
    @{Row=51; D="0.664"; E="  +3.71%  "}
)

foreach ($item in $priceVolumeUpdates) {
    if ($item.D -ne $null) {
        $dCell = $ws.Cells.Item($item.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# Rows whose ranking moved, so the Coin/Link were replaced along with new Price/Volume(1h).
$fullRowUpdates = @(
    @{Row=36; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="7.25"; E="  +24.07%  "},
    @{Row=37; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.148"; E="  -3.31%  "},
    @{Row=45; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="157.77"; E="  -0.65%  "},
    @{Row=46; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="188.89"; E="  -0.92%  "},
    @{Row=47; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.43"; E="  +5.60%  "},
    @{Row=48; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="46.75"; E="  +8.66%  "},
    @{Row=49; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.48"; E="  +3.48%  "},
    @{Row=50; B="Mantle"; C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D="0.789"; E="  -0.51%  "}
)

foreach ($item in $fullRowUpdates) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $dCell.Style = "Normal"
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}
